$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Brighton v Leeds'
$ws.Range("B2").Value = 'Brighton'
$ws.Range("C2").Value = 'England Premier League'
$ws.Range("D2").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E2").Value = '112/141 Win Tips'
$ws.Range("F2").Value = "'79"
$ws.Range("G2").Value = "'2.00"

# Row 3
$ws.Range("A3").Value = 'Fulham v Wolverhampton'
$ws.Range("B3").Value = 'Fulham'
$ws.Range("C3").Value = 'England Premier League'
$ws.Range("D3").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E3").Value = '98/137 Win Tips'
$ws.Range("F3").Value = "'72"
$ws.Range("G3").Value = "'1.80"

# Row 4
$ws.Range("A4").Value = 'Nottm Forest v Man Utd'
$ws.Range("B4").Value = 'Man Utd'
$ws.Range("C4").Value = 'England Premier League'
$ws.Range("D4").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E4").Value = '92/128 Win Tips'
$ws.Range("F4").Value = "'72"
$ws.Range("G4").Value = "'2.10"

# Row 5
$ws.Range("A5").Value = 'Burnley v Arsenal'
$ws.Range("B5").Value = 'Arsenal'
$ws.Range("C5").Value = 'England Premier League'
$ws.Range("D5").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E5").Value = '90/121 Win Tips'
$ws.Range("F5").Value = "'74"
$ws.Range("G5").Value = "'1.25"

# Row 6
$ws.Range("A6").Value = 'Liverpool v Aston Villa'
$ws.Range("B6").Value = 'Liverpool'
$ws.Range("C6").Value = 'England Premier League'
$ws.Range("D6").Value = '2025-11-01T20:00:00.000Z'
$ws.Range("E6").Value = '53/97 Win Tips'
$ws.Range("F6").Value = "'55"
$ws.Range("G6").Value = "'1.70"

# Row 7
$ws.Range("A7").Value = 'Crystal Palace v Brentford'
$ws.Range("B7").Value = 'Draw'
$ws.Range("C7").Value = 'England Premier League'
$ws.Range("D7").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E7").Value = '43/107 Win Tips'
$ws.Range("F7").Value = "'40"
$ws.Range("G7").Value = "'3.60"

# Row 8
$ws.Range("A8").Value = 'Tottenham v Chelsea'
$ws.Range("B8").Value = 'Chelsea'
$ws.Range("C8").Value = 'England Premier League'
$ws.Range("D8").Value = '2025-11-01T17:30:00.000Z'
$ws.Range("E8").Value = '40/92 Win Tips'
$ws.Range("F8").Value = "'43"
$ws.Range("G8").Value = "'2.50"

# Row 9
$ws.Range("A9").Value = 'West Ham v Newcastle'
$ws.Range("B9").Value = 'Newcastle'
$ws.Range("C9").Value = 'England Premier League'
$ws.Range("D9").Value = '2025-11-02T14:00:00.000Z'
$ws.Range("E9").Value = '40/51 Win Tips'
$ws.Range("F9").Value = "'78"
$ws.Range("G9").Value = "'1.67"

# Row 10
$ws.Range("A10").Value = 'Man City v Bournemouth'
$ws.Range("B10").Value = 'Man City'
$ws.Range("C10").Value = 'England Premier League'
$ws.Range("D10").Value = '2025-11-02T16:30:00.000Z'
$ws.Range("E10").Value = '39/49 Win Tips'
$ws.Range("F10").Value = "'80"
$ws.Range("G10").Value = "'1.53"

# Row 11
$ws.Range("A11").Value = 'Real Madrid v Valencia'
$ws.Range("B11").Value = 'Real Madrid'
$ws.Range("C11").Value = 'Spain Primera Liga'
$ws.Range("D11").Value = '2025-11-01T20:00:00.000Z'
$ws.Range("E11").Value = '26/30 Win Tips'
$ws.Range("F11").Value = "'87"
$ws.Range("G11").Value = "'1.18"

# Row 12
$ws.Range("A12").Value = 'Heidenheim v Eintracht Frankfurt'
$ws.Range("B12").Value = 'Eintracht Frankfurt'
$ws.Range("C12").Value = 'Germany Bundesliga I'
$ws.Range("D12").Value = '2025-11-01T14:30:00.000Z'
$ws.Range("E12").Value = '23/23 Win Tips'
$ws.Range("F12").Value = "'100"
$ws.Range("G12").Value = "'1.91"

# Row 13
$ws.Range("A13").Value = 'PSG v Nice'
$ws.Range("B13").Value = 'PSG'
$ws.Range("C13").Value = 'France Ligue 1'
$ws.Range("D13").Value = '2025-11-01T16:00:00.000Z'
$ws.Range("E13").Value = '23/25 Win Tips'
$ws.Range("F13").Value = "'92"
$ws.Range("G13").Value = "'1.20"

# Row 14
$ws.Range("A14").Value = 'Charlton v Swansea'
$ws.Range("B14").Value = 'Charlton'
$ws.Range("C14").Value = 'England Championship'
$ws.Range("D14").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E14").Value = '21/30 Win Tips'
$ws.Range("F14").Value = "'70"
$ws.Range("G14").Value = "'2.15"

# Row 15
$ws.Range("A15").Value = 'Bayern Munich v Bayer Leverkusen'
$ws.Range("B15").Value = 'Bayern Munich'
$ws.Range("C15").Value = 'Germany Bundesliga I'
$ws.Range("D15").Value = '2025-11-01T17:30:00.000Z'
$ws.Range("E15").Value = '21/24 Win Tips'
$ws.Range("F15").Value = "'88"
$ws.Range("G15").Value = "'1.22"

# Row 16
$ws.Range("A16").Value = 'Sunderland v Everton'
$ws.Range("B16").Value = 'Sunderland'
$ws.Range("C16").Value = 'England Premier League'
$ws.Range("D16").Value = '2025-11-03T20:00:00.000Z'
$ws.Range("E16").Value = '21/32 Win Tips'
$ws.Range("F16").Value = "'66"
$ws.Range("G16").Value = "'2.88"

# Row 17
$ws.Range("A17").Value = 'Napoli v Como'
$ws.Range("B17").Value = 'Napoli'
$ws.Range("C17").Value = 'Italy Serie A'
$ws.Range("D17").Value = '2025-11-01T17:00:00.000Z'
$ws.Range("E17").Value = '20/23 Win Tips'
$ws.Range("F17").Value = "'87"
$ws.Range("G17").Value = "'1.95"

# Row 18
$ws.Range("A18").Value = 'Oxford Utd v Millwall'
$ws.Range("B18").Value = 'Millwall'
$ws.Range("C18").Value = 'England Championship'
$ws.Range("D18").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E18").Value = '19/26 Win Tips'
$ws.Range("F18").Value = "'73"
$ws.Range("G18").Value = "'2.85"

# Row 19
$ws.Range("A19").Value = 'Atletico Madrid v Sevilla'
$ws.Range("B19").Value = 'Atletico Madrid'
$ws.Range("C19").Value = 'Spain Primera Liga'
$ws.Range("D19").Value = '2025-11-01T15:15:00.000Z'
$ws.Range("E19").Value = '19/22 Win Tips'
$ws.Range("F19").Value = "'86"
$ws.Range("G19").Value = "'1.35"

# Row 20
$ws.Range("A20").Value = 'Sheff Utd v Derby'
$ws.Range("B20").Value = 'Sheff Utd'
$ws.Range("C20").Value = 'England Championship'
$ws.Range("D20").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E20").Value = '17/27 Win Tips'
$ws.Range("F20").Value = "'63"
$ws.Range("G20").Value = "'1.80"

# Row 21
$ws.Range("A21").Value = 'Birmingham v Portsmouth'
$ws.Range("B21").Value = 'Birmingham'
$ws.Range("C21").Value = 'England Championship'
$ws.Range("D21").Value = '2025-11-01T15:00:00.000Z'
$ws.Range("E21").Value = '15/26 Win Tips'
$ws.Range("F21").Value = "'58"
$ws.Range("G21").Value = "'1.75"

# Row 22
$ws.Range("B22").Value = 'PSG'
$ws.Range("G22").Value = "'5.50"

# Row 23
$ws.Range("E23").Value = '4/4 Win Tips'
$ws.Range("F23").ClearContents()

# The leading apostrophes above make Excel apply a "quote prefix" / text
# number format to these cells; reset them back to the default General
# style so the text values are stored the same way as the rest of the sheet.
$ws.Range("F2:G21").Style = "Normal"
$ws.Range("G22").Style = "Normal"
